# Emmersive localization workbook: add "active blacklist/whitelist character" strings
# feat(EM): add blacklist mode

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# New shared strings must be appended in this exact order to match the
# original authoring order in the workbook's shared string table:
#   em_ui_active_blacklist, em_ui_active_whitelist,
#   活跃的黑名单角色, 活跃的白名单角色,
#   アクティブなブラックリスト, アクティブなホワイトリスト
$ws.Range("A116").Value = "em_ui_active_blacklist"
$ws.Range("A117").Value = "em_ui_active_whitelist"
$ws.Range("D116").Value = "活跃的黑名单角色"
$ws.Range("D117").Value = "活跃的白名单角色"
$ws.Range("C116").Value = "アクティブなブラックリスト"
$ws.Range("C117").Value = "アクティブなホワイトリスト"

# Reflect the author's final view/selection state (scroll position + selection)
$win = $excel.ActiveWindow
try {
    $win.ScrollRow = 102
    $win.ScrollColumn = 1
} catch {}

$ws.Range("A116:D117").Select()
